$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-30 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-01 Thursday", 2) | Out-Null
$d.Content.Find.Execute("164÷7=23, 3", $true, $false, $false, $false, $false, $true, 1, $false, "152÷2=76, 0", 2) | Out-Null
$d.Content.Find.Execute("890÷6=148, 2", $true, $false, $false, $false, $false, $true, 1, $false, "342÷9=38, 0", 2) | Out-Null
$d.Content.Find.Execute("492÷6=82, 0", $true, $false, $false, $false, $false, $true, 1, $false, "775÷2=387, 1", 2) | Out-Null
$d.Content.Find.Execute("565÷9=62, 7", $true, $false, $false, $false, $false, $true, 1, $false, "928÷7=132, 4", 2) | Out-Null
$d.Content.Find.Execute("423÷6=70, 3", $true, $false, $false, $false, $false, $true, 1, $false, "874÷8=109, 2", 2) | Out-Null
$d.Content.Find.Execute("340÷4=85, 0", $true, $false, $false, $false, $false, $true, 1, $false, "713÷9=79, 2", 2) | Out-Null
$d.Content.Find.Execute("922÷7=131, 5", $true, $false, $false, $false, $false, $true, 1, $false, "712÷3=237, 1", 2) | Out-Null
$d.Content.Find.Execute("996÷4=249, 0", $true, $false, $false, $false, $false, $true, 1, $false, "542÷6=90, 2", 2) | Out-Null
$d.Content.Find.Execute("100÷9=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "108÷2=54, 0", 2) | Out-Null
$d.Content.Find.Execute("225÷4=56, 1", $true, $false, $false, $false, $false, $true, 1, $false, "896÷9=99, 5", 2) | Out-Null
$d.Content.Find.Execute("188÷8=23, 4", $true, $false, $false, $false, $false, $true, 1, $false, "334÷3=111, 1", 2) | Out-Null
$d.Content.Find.Execute("394÷5=78, 4", $true, $false, $false, $false, $false, $true, 1, $false, "587÷8=73, 3", 2) | Out-Null
$d.Content.Find.Execute("485÷9=53, 8", $true, $false, $false, $false, $false, $true, 1, $false, "540÷7=77, 1", 2) | Out-Null
$d.Content.Find.Execute("651÷4=162, 3", $true, $false, $false, $false, $false, $true, 1, $false, "721÷6=120, 1", 2) | Out-Null
$d.Content.Find.Execute("275÷3=91, 2", $true, $false, $false, $false, $false, $true, 1, $false, "742÷4=185, 2", 2) | Out-Null
$d.Content.Find.Execute("465÷6=77, 3", $true, $false, $false, $false, $false, $true, 1, $false, "180÷5=36, 0", 2) | Out-Null
$d.Content.Find.Execute("234÷2=117, 0", $true, $false, $false, $false, $false, $true, 1, $false, "587÷4=146, 3", 2) | Out-Null
$d.Content.Find.Execute("571÷9=63, 4", $true, $false, $false, $false, $false, $true, 1, $false, "342÷4=85, 2", 2) | Out-Null
$d.Content.Find.Execute("624÷9=69, 3", $true, $false, $false, $false, $false, $true, 1, $false, "478÷9=53, 1", 2) | Out-Null
$d.Content.Find.Execute("731÷2=365, 1", $true, $false, $false, $false, $false, $true, 1, $false, "248÷4=62, 0", 2) | Out-Null
$d.Content.Find.Execute("641÷5=128, 1", $true, $false, $false, $false, $false, $true, 1, $false, "403÷7=57, 4", 2) | Out-Null
$d.Content.Find.Execute("458÷2=229, 0", $true, $false, $false, $false, $false, $true, 1, $false, "766÷2=383, 0", 2) | Out-Null
$d.Content.Find.Execute("793÷5=158, 3", $true, $false, $false, $false, $false, $true, 1, $false, "550÷9=61, 1", 2) | Out-Null
$d.Content.Find.Execute("332÷6=55, 2", $true, $false, $false, $false, $false, $true, 1, $false, "298÷6=49, 4", 2) | Out-Null
$d.Content.Find.Execute("363÷2=181, 1", $true, $false, $false, $false, $false, $true, 1, $false, "534÷7=76, 2", 2) | Out-Null
